# Applies the "Add data for 2022-07-18" update to cta-violent-crime-ytd.xlsx
# This updates year-to-date violent-crime counts (through July 18) across the
# Citywide Totals sheet, the By Neighborhood summary sheet, and several
# individual neighborhood sheets that received new/updated daily data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Citywide Totals
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 51
$ws.Range("H2").Value = 61
$ws.Range("F3").Value = 78
$ws.Range("I3").Value = 105
$ws.Range("B6").Value = 211
$ws.Range("C6").Value = 262
$ws.Range("F6").Value = 300
$ws.Range("G6").Value = 267
$ws.Range("H6").Value = 235
$ws.Range("I6").Value = 298
$ws.Range("B7").Value = 284
$ws.Range("C7").Value = 353
$ws.Range("D7").Value = 377
$ws.Range("F7").Value = 430
$ws.Range("G7").Value = 393
$ws.Range("H7").Value = 366
$ws.Range("I7").Value = 482

# ---------------------------------------------------------------------------
# By Neighborhood
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("C18").Value = 8
$ws.Range("F19").Value = 4
$ws.Range("B27").Value = 24
$ws.Range("H28").Value = 4
$ws.Range("I42").Value = 3
$ws.Range("F49").Value = 17
$ws.Range("I49").Value = 8
$ws.Range("D52").Value = 53
$ws.Range("F52").Value = 41
$ws.Range("H52").Value = 41
$ws.Range("I52").Value = 75
$ws.Range("I60").Value = 2
$ws.Range("C64").Value = 10
$ws.Range("G73").Value = 8
$ws.Range("F90").Value = 6
$ws.Range("B97").Value = 284
$ws.Range("C97").Value = 353
$ws.Range("D97").Value = 377
$ws.Range("F97").Value = 430
$ws.Range("G97").Value = 393
$ws.Range("H97").Value = 366
$ws.Range("I97").Value = 482

# ---------------------------------------------------------------------------
# Chicago Lawn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 4

# ---------------------------------------------------------------------------
# Little Italy, UIC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I3").Value = 1
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 17
$ws.Range("I6").Value = 8

# ---------------------------------------------------------------------------
# Englewood
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B5").Value = 24
$ws.Range("B6").Value = 24

# ---------------------------------------------------------------------------
# Loop
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D2").Value = 9
$ws.Range("I3").Value = 13
$ws.Range("F6").Value = 32
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 50
$ws.Range("D7").Value = 53
$ws.Range("F7").Value = 41
$ws.Range("H7").Value = 41
$ws.Range("I7").Value = 75

# ---------------------------------------------------------------------------
# West Loop
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 6

# ---------------------------------------------------------------------------
# North Lawndale
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("C5").Value = 9
$ws.Range("C6").Value = 10

# ---------------------------------------------------------------------------
# River North
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("River North")
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 8

# ---------------------------------------------------------------------------
# Fuller Park
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("H2").Value = 1
$ws.Range("H6").Value = 4

# ---------------------------------------------------------------------------
# Chatham
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 8

# ---------------------------------------------------------------------------
# Irving Park
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 3
